$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.981.65"
$ws.Range("E2").Value = "  +3.09%  "

$ws.Range("D3").Value = "2.318.15"
$ws.Range("E3").Value = "  +1.14%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.45"
$ws.Range("E5").Value = "  +1.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.43"
$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.581"
$ws.Range("E8").Value = "  -0.39%  "

$ws.Range("D9").Value = "2.316.01"
$ws.Range("E9").Value = "  +1.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.100"
$ws.Range("E10").Value = "  +0.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.49"
$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("E12").Value = "  +0.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.54"
$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("D15").Value = "2.734.28"
$ws.Range("E15").Value = "  +1.22%  "

$ws.Range("D16").Value = "59.964.41"
$ws.Range("E16").Value = "  +3.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").Value = "2.313.76"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.53"
$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.13"
$ws.Range("E20").Value = "  -1.47%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.26"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.63"
$ws.Range("E22").Value = "  +3.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.995"
$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.56"
$ws.Range("E24").Value = "  +1.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.79"
$ws.Range("E27").Value = "  -2.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.33"
$ws.Range("E28").Value = "  +5.08%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.96"
$ws.Range("E29").Value = "  +2.00%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.73"
$ws.Range("E30").Value = "  +1.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  +11.06%  "

$ws.Range("D32").Value = "0.0₃0722"
$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.88"
$ws.Range("E33").Value = "  +2.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.37"
$ws.Range("E34").Value = "  +11.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.379"
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.89"
$ws.Range("E37").Value = "  +0.60%  "

$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.04"
$ws.Range("E39").Value = "  +4.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "320.83"
$ws.Range("E40").Value = "  +11.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.92"
$ws.Range("E41").Value = "  -1.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.50"
$ws.Range("E42").Value = "  +0.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.01"
$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.43"
$ws.Range("E44").Value = "  +0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0942"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.36"
$ws.Range("E46").Value = "  +7.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0495"
$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.556"
$ws.Range("E48").Value = "  +0.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0211"
$ws.Range("E49").Value = "  +0.84%  "

$ws.Range("D50").Value = "0.0₆0212"
$ws.Range("E50").Value = "  +12.01%  "

$ws.Range("E51").Value = "  +0.78%  "

